# Update "想去人数" (F column) counts on several rows across sheets
# 展览 (sheet "展览")
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 329
$ws1.Range("F5").Value = 179
$ws1.Range("F6").Value = 679
$ws1.Range("F8").Value = 487
$ws1.Range("F10").Value = 529
$ws1.Range("F11").Value = 412
$ws1.Range("F14").Value = 118
$ws1.Range("F15").Value = 202

# 本地生活 (sheet "本地生活")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6227
$ws3.Range("F4").Value = 756
$ws3.Range("F5").Value = 1828

# 全部类型 (sheet "全部类型")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6227
$ws4.Range("F4").Value = 756
$ws4.Range("F5").Value = 1828
$ws4.Range("F6").Value = 329
$ws4.Range("F12").Value = 179
$ws4.Range("F15").Value = 679
$ws4.Range("F19").Value = 487
$ws4.Range("F22").Value = 529
$ws4.Range("F24").Value = 412
$ws4.Range("F29").Value = 118
$ws4.Range("F35").Value = 202
